# Update existing row 2 and append two new data rows to the "HurtoListado"
# style report sheet, matching the same layout/format already used for row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: replace with the new incident record ------------------------
$ws.Range("A2").Value = "DTSC"
$ws.Range("B2").Value = "2020-10-06"
$ws.Range("C2").Value = "Llamada de amenaza de Bomba en el CAT Aguilera No. 401"
$ws.Range("D2").Value = "Santiago de Cuba"
$ws.Range("E2").Value = "17488/20"
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = "DTSC-10-20-0106"

# --- Row 3: new incident record ------------------------------------------
$ws.Range("A3").Value = "DVLH"
$ws.Range("B3").Value = "2020-12-02"
$ws.Range("C3").Value = "Fractura de parabrisa delantero de un auto, por agresión de una ciudadana Calle Águila, entre Dragones y Barcelona"
$ws.Range("D3").Value = "Centro Habana"
$ws.Range("E3").Value = "65472/20"
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = "DVLH-12-20-0228"

# --- Row 4: new incident record (last row - closes the table border) -----
$ws.Range("A4").Value = "DTSR"
$ws.Range("B4").Value = "2020-12-10"
$ws.Range("C4").Value = "Individuo que fractura el cristal de la puerta de entrada Calle 10 de Octubre. No. 1251,  entre Carmen y Vista Alegre"
$ws.Range("D4").Value = "Díez de Octubre"
$ws.Range("E4").Value = "67102/20"
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = "DVLH-12-20-0234"

# --- Formatting: copy row 2's style down to rows 3 and 4, then apply a
# bottom border across row 4 (A:K) since it is now the last row of the table.
$ws.Range("A2:K2").Copy() | Out-Null
$ws.Range("A3:K3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A4:K4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A4:K4").Borders.Item(9).LineStyle = 1     # xlEdgeBottom, xlContinuous
$ws.Range("A4:K4").Borders.Item(9).Weight = 2        # xlThin

$excel.CutCopyMode = 0
